# Update instrument reference designators from GA05MOAS-GL003 to GA05MOAS-GL496
$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCalInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# Moorings sheet: Ref Des in A2
$wsMoorings.Range("A2").Value = "GA05MOAS-GL496"

# Asset_Cal_Info sheet: Ref Des values in column A
$wsCalInfo.Range("A2").Value  = "GA05MOAS-GL496-01-FLORDM000"
$wsCalInfo.Range("A3").Value  = "GA05MOAS-GL496-01-FLORDM000"
$wsCalInfo.Range("A4").Value  = "GA05MOAS-GL496-01-FLORDM000"
$wsCalInfo.Range("A5").Value  = "GA05MOAS-GL496-01-FLORDM000"
$wsCalInfo.Range("A7").Value  = "GA05MOAS-GL496-02-DOSTAM000"
$wsCalInfo.Range("A9").Value  = "GA05MOAS-GL496-04-CTDGVM000"
$wsCalInfo.Range("A11").Value = "GA05MOAS-GL496-00-ENG000000"

# Update the selected/active cell on the Moorings sheet view
$wsMoorings.Activate()
$wsMoorings.Range("E19").Select()
